$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The I_holdings export had merged in OHLC + ticker data scraped from other
# tickers' rows (WIX, NOC, CDNS, ...). Restore each dated row to PRO's own
# open/close/high/low price, shares_outstanding, and ticker label.
$priceData = @(
    @{Row=2; Open=24.45999908447266; Close=22.22999954223633; High=27.22999954223633; Low=20.8700008392334}
    @{Row=3; Open=21.31999969482422; Close=21.8700008392334; High=21.92000007629395; Low=20.10000038146973}
    @{Row=4; Open=22.01000022888184; Close=24.02000045776367; High=25; Low=21.56999969482422}
    @{Row=5; Open=22.54000091552734; Close=12.27999973297119; High=22.63999938964844; Low=11.64000034332275}
    @{Row=6; Open=11.67000007629394; Close=11.76000022888184; High=11.9399995803833; Low=10.60999965667725}
    @{Row=7; Open=17.3799991607666; Close=18.57999992370605; High=18.71999931335449; Low=16.76000022888184}
    @{Row=8; Open=22.53000068664551; Close=21.92000007629395; High=22.8799991607666; Low=20.55999946594238}
    @{Row=9; Open=21.65999984741211; Close=22.3799991607666; High=22.42000007629395; Low=20.34000015258789}
    @{Row=10; Open=24.1200008392334; Close=24.64999961853028; High=24.84000015258789; Low=22.3799991607666}
    @{Row=11; Open=27.59000015258789; Close=28.82999992370605; High=29.69000053405762; Low=26.96999931335449}
    @{Row=12; Open=24.20000076293945; Close=22.59000015258789; High=25.80999946594238; Low=21.36000061035156}
    @{Row=13; Open=26.45000076293945; Close=29.04000091552734; High=29.78000068664551; Low=26.22999954223633}
    @{Row=14; Open=33.02999877929688; Close=29.52000045776367; High=34.13000106811523; Low=29.19000053405762}
    @{Row=15; Open=37.08000183105469; Close=37.13999938964844; High=41.09999847412109; Low=36.40000152587891}
    @{Row=16; Open=35.15000152587891; Close=32.91999816894531; High=35.56000137329102; Low=30.56999969482422}
    @{Row=17; Open=30.89999961853028; Close=34.59999847412109; High=34.9900016784668; Low=29.36000061035156}
    @{Row=18; Open=41.97000122070312; Close=51.2400016784668; High=53.34999847412109; Low=41.54999923706055}
    @{Row=19; Open=63.95000076293945; Close=72.36000061035156; High=74.7300033569336; Low=62.36000061035156}
    @{Row=20; Open=59.88000106811523; Close=51.2400016784668; High=61.68000030517578; Low=45.47999954223633}
    @{Row=21; Open=60.5099983215332; Close=60; High=68.80999755859375; Low=59.54000091552734}
    @{Row=22; Open=29.68000030517578; Close=34.38999938964844; High=34.79999923706055; Low=24.56999969482422}
    @{Row=23; Open=44.20999908447266; Close=32.63000106811523; High=46.02000045776367; Low=30.30999946594238}
    @{Row=24; Open=32.15999984741211; Close=28.17000007629395; High=35.84999847412109; Low=24.54999923706055}
    @{Row=25; Open=50.7400016784668; Close=42.13999938964844; High=50.95000076293945; Low=41.02000045776367}
    @{Row=26; Open=43.06000137329102; Close=42.97999954223633; High=47.72999954223633; Low=42.27000045776367}
    @{Row=27; Open=45.7599983215332; Close=43.41999816894531; High=47.61999893188477; Low=41.72000122070312}
    @{Row=28; Open=35.81999969482422; Close=30; High=37.11000061035156; Low=29.71999931335449}
    @{Row=29; Open=34.72999954223633; Close=27.70999908447266; High=36.29000091552734; Low=25.04999923706055}
    @{Row=30; Open=33.31999969482422; Close=27.93000030517578; High=35.11999893188477; Low=27.28000068664551}
    @{Row=31; Open=26.03000068664551; Close=24.3700008392334; High=29.02000045776367; Low=22.20999908447266}
    @{Row=32; Open=25.15999984741211; Close=24.95000076293945; High=28.09000015258789; Low=23.85000038146973}
    @{Row=33; Open=24.71999931335449; Close=25.20000076293945; High=26.18000030517578; Low=22.76000022888184}
    @{Row=34; Open=27.09000015258789; Close=28.3700008392334; High=30.20000076293945; Low=25.52000045776367}
    @{Row=35; Open=30.64999961853028; Close=38; High=38.95999908447266; Low=29.47999954223633}
    @{Row=36; Open=34.59999847412109; Close=31.14999961853028; High=35.31999969482422; Low=30.95000076293945}
    @{Row=37; Open=38.25; Close=34.41999816894531; High=38.25; Low=32.84999847412109}
    @{Row=38; Open=36.38000106811523; Close=32.75; High=38.70999908447266; Low=32.72000122070312}
    @{Row=39; Open=28.75; Close=24.10000038146973; High=28.75; Low=23.29000091552734}
    @{Row=40; Open=18.46999931335449; Close=19.79999923706055; High=23; Low=17.48999977111816}
    @{Row=41; Open=22.34000015258789; Close=23.6200008392334; High=25.13999938964844; Low=21}
    @{Row=42; Open=18.98999977111816; Close=17.07999992370605; High=19.54999923706055; Low=14.8100004196167}
    @{Row=43; Open=15.68000030517578; Close=15.6899995803833; High=17.10000038146973; Low=14.75}
)

foreach ($rec in $priceData) {
    $r = $rec.Row
    $ws.Cells.Item($r, 4).Value = $rec.Open
    $ws.Cells.Item($r, 5).Value = $rec.Close
    $ws.Cells.Item($r, 6).Value = $rec.High
    $ws.Cells.Item($r, 7).Value = $rec.Low
    $ws.Cells.Item($r, 8).Value = 48253392
    $ws.Cells.Item($r, 9).Value = "PRO"
}
